$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9740301966667175
$ws.Range("B1").Value = 1.965752124786377
$ws.Range("C1").Value = 2.928719282150269
$ws.Range("D1").Value = 2.378083229064941
$ws.Range("E1").Value = 0.8414344787597656
